$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells keep their literal text (avoid auto-numeric conversion to a Double,
# which would silently drop significant trailing zeros, e.g. "113.80" -> 113.8)
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply updated coin data (values scraped fresh by the GitHub Actions job)
$ws.Range('D2').Value = '26.646.95'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').Value = '1.826.85'
$ws.Range('E3').Value = '  +1.78%  '
$ws.Range('E4').Value = '  +0.41%  '
$ws.Range('D5').Value = '309.36'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('D7').Value = '0.4675'
$ws.Range('E7').Value = '  +3.50%  '
$ws.Range('D8').Value = '0.3596'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  +0.80%  '
$ws.Range('D10').Value = '0.9042'
$ws.Range('E10').Value = '  +2.17%  '
$ws.Range('D11').Value = '0.07684'
$ws.Range('E11').Value = '  -0.81%  '
$ws.Range('D12').Value = '19.43'
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').Value = '1.849.89'
$ws.Range('E13').Value = '  +3.38%  '
$ws.Range('D14').Value = '5.266'
$ws.Range('E14').Value = '  -0.32%  '
$ws.Range('D15').Value = '6.377'
$ws.Range('E15').Value = '  +0.92%  '
$ws.Range('D16').Value = '87.78'
$ws.Range('E16').Value = '  +3.37%  '
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').Value = '0.000008558'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('D20').Value = '26.661.94'
$ws.Range('E20').Value = '  +1.08%  '
$ws.Range('E21').Value = '  -0.30%  '
$ws.Range('D22').Value = '5.026'
$ws.Range('E22').Value = '  +1.12%  '
$ws.Range('D23').Value = '10.56'
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('D24').Value = '1.907'
$ws.Range('E24').Value = '  -3.29%  '
$ws.Range('D25').Value = '152.92'
$ws.Range('E25').Value = '  +1.11%  '
$ws.Range('D26').Value = '17.94'
$ws.Range('E26').Value = '  +0.59%  '
$ws.Range('D27').Value = '1.998'
$ws.Range('E27').Value = '  -1.59%  '
$ws.Range('D28').Value = '113.80'
$ws.Range('E28').Value = '  +1.65%  '
$ws.Range('D29').Value = '4.870'
$ws.Range('E29').Value = '  +0.54%  '
$ws.Range('D30').Value = '0.08825'
$ws.Range('E30').Value = '  +1.60%  '
$ws.Range('D31').Value = '3.148'
$ws.Range('E31').Value = '  +2.10%  '
$ws.Range('D32').Value = '2.859'
$ws.Range('E32').Value = '  +4.22%  '
$ws.Range('D33').Value = '1.170'
$ws.Range('E33').Value = '  +5.98%  '
$ws.Range('D34').Value = '0.7375'
$ws.Range('E34').Value = '  +2.24%  '
$ws.Range('D35').Value = '4.434'
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('D36').Value = '1.077'
$ws.Range('E36').Value = '  +0.87%  '
$ws.Range('D37').Value = '0.01934'
$ws.Range('E37').Value = '  +0.19%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '2.944'
$ws.Range('E38').Value = '  +3.09%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.05159'
$ws.Range('E39').Value = '  +1.42%  '
$ws.Range('D40').Value = '6.878'
$ws.Range('E40').Value = '  +0.76%  '
$ws.Range('D41').Value = '0.5064'
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('E42').Value = '  -1.11%  '
$ws.Range('D43').Value = '8.059'
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('D44').Value = '1.008'
$ws.Range('E44').Value = '  +0.36%  '
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('D46').Value = '10.02'
$ws.Range('E46').Value = '  +1.61%  '
$ws.Range('D47').Value = '98.45'
$ws.Range('E47').Value = '  -2.17%  '
$ws.Range('D48').Value = '1.574'
$ws.Range('E48').Value = '  +0.59%  '
$ws.Range('D49').Value = '0.06026'
$ws.Range('E49').Value = '  +1.14%  '
$ws.Range('D50').Value = '64.02'
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('D51').Value = '35.83'
$ws.Range('E51').Value = '  -0.67%  '
